# Import location (suburb/postal code/city) and a stable currency-style
# numeric id onto the sales report's Table1, replacing the old
# "Days/Value/Active" sample columns with a single location record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename the table headers (renaming the header cell also renames the
#     bound ListColumn) -------------------------------------------------
$ws.Range("B1").Value = "SuburbName"
$ws.Range("C1").Value = "PostalCode"
$ws.Range("D1").Value = "CityId"

# --- Replace the single remaining data row with the imported record ----
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "Hattflied"
$ws.Range("C2").Value = 231
$ws.Range("D2").Value = 1

# --- Drop the old sample rows (3-10) ------------------------------------
$ws.Range("A3:D10").ClearContents()

# --- Shrink the table/autofilter range down to the remaining data ------
$tbl = $ws.ListObjects.Item(1)
$tbl.Resize($ws.Range("A1:D2"))

# --- Page setup was touched as part of the re-import (portrait print) --
$ws.PageSetup.Orientation = 1

# --- Restore the selection left behind by the edit ----------------------
$ws.Range("E12").Select()
